$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the Amazon Order 9 refund (row 22, col E):
#    =1537.58-8.18  ->  =1537.58-8.18-326
# ------------------------------------------------------------------
$ws.Range("E22").Formula = "=1537.58-8.18-326"

# ------------------------------------------------------------------
# 2. Insert three new rows (25, 26, 27) above the old "Total Expenses"
#    row, pushing it (and the two footnote rows after it) down by three.
# ------------------------------------------------------------------
$ws.Range("A25:F27").Insert()

# Row 25: Arrow Order 5 (no receipt link; just a blank hyperlink-styled cell)
$ws.Range("A25").Value = "Arrow Order 5"
$ws.Range("B25").Value = "2/24/2018"
$ws.Range("C25").Value = "Morgan"
$ws.Range("D25").Style = "Hyperlink"
$ws.Range("E25").Value = 58.73
$ws.Range("F25").Value = "N-Channel MOSFETs"

# Row 26: Amazon Order 11
$ws.Range("A26").Value = "Amazon Order 11"
$ws.Range("B26").Value = "2/24/2018"
$ws.Range("C26").Value = "Morgan"
$ws.Range("D26").Value = "Amazon Order 11.pdf"
$ws.Hyperlinks.Add($ws.Range("D26"), "Amazon%20Order%2011.pdf") | Out-Null
$ws.Range("E26").Value = 209.04

# ------------------------------------------------------------------
# 3. Update the refund footnote text (now row 30) to the longer
#    version that also mentions the returned USB cables.
# ------------------------------------------------------------------
$ws.Range("D30").Value = "***Refunded `$8.18 for an incorrect shipment and an additional `$326.00 for returned USB cables"

# Row 27: JLCPCB Order
$ws.Range("A27").Value = "JLCPCB Order"
$ws.Range("B27").Value = "2/24/2018"
$ws.Range("C27").Value = "Morgan"
$ws.Range("D27").Value = "JLCPCB Order.pdf"
$ws.Hyperlinks.Add($ws.Range("D27"), "JLCPCB Order.pdf") | Out-Null
$ws.Range("E27").Value = 58.07

# Match date-format / currency-format styling used elsewhere on the sheet
# for these new data rows.
$ws.Range("B25:B27").NumberFormat = "mm/dd/yy;@"
$ws.Range("E25:E27").NumberFormat = "$#,##0.00"

# ------------------------------------------------------------------
# 4. Fix up the "Total Expenses" row (now row 28) so the SUM range
#    covers the newly inserted rows.
# ------------------------------------------------------------------
$ws.Range("E28").Formula = "=SUM(E2:E27)"

# ------------------------------------------------------------------
# 5. Leave the selection near the bottom of the new data, matching
#    where the author was working when they saved the file.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A28").Select()

Write-Output "done"
